# Fruta / hortaliza, semanal
# Insert a new data row at row 56 (pushing the existing rows 56-90 down to
# 57-91) and populate the new row with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 56:90 down to 57:91, leaving an empty row 56 behind.
$ws.Rows.Item(56).Insert()

# Fill in the newly inserted row 56 with the new weekly record.
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = 44574
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 100112032
$ws.Range("G56").Value = "Zapallo italiano"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 100
$ws.Range("K56").Value = 13000
$ws.Range("L56").Value = 14000
$ws.Range("M56").Value = 13500
$ws.Range("N56").Value = '$/caja 60 unidades'
$ws.Range("O56").Value = "Región de O'Higgins"
$ws.Range("P56").Value = 225
$ws.Range("Q56").Value = 60
$ws.Range("R56").Value = "Hortaliza"
